$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '29.507.58'
$ws.Range("E2").Value = '  -0.86%  '
$ws.Range("D3").Value = '1.850.43'
$ws.Range("E3").Value = '  -0.44%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.9996'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  +0.08%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '241.82'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -1.10%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.6253'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -2.71%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.000'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +0.11%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '47.94'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +1.07%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.07559'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +0.35%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.2977'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -0.28%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '24.27'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -1.02%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.07675'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +0.06%  '
$ws.Range("D13").Value = '1.907.11'
$ws.Range("E13").Value = '  +2.60%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.015'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -0.50%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.6853'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -1.04%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '83.82'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -0.13%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.000009738'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -2.14%  '
$ws.Range("D18").Value = '2.145.32'
$ws.Range("E18").Value = '  +1.31%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '6.218'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +1.84%  '
$ws.Range("D20").Value = '29.569.71'
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '234.48'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -0.85%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '12.49'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -1.56%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '1.0000'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +0.05%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '7.607'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +1.19%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '1.001'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +0.08%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '155.79'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -2.00%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.1389'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -2.42%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '8.423'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -1.28%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '17.72'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -1.05%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.481'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -0.77%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.05829'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -6.17%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.261'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -2.01%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '4.106'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -1.39%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '4.033'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -1.72%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.892'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -0.18%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.171'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -0.05%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.7185'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -1.56%  '
$ws.Range("E38").Value = '  -0.60%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.803'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -0.77%  '
$ws.Range("D40").Value = '1.237.02'
$ws.Range("E40").Value = '  +2.99%  '
$ws.Range("E41").Value = '  -0.57%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.9131'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -1.04%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '6.132'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -1.88%  '
$ws.Range("D44").Value = '2.051.24'
$ws.Range("E44").Value = '  +1.14%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.9995'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -0.03%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '102.84'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +0.80%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '67.42'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +1.34%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '7.322'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +9.50%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '9.161'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -0.29%  '
$ws.Range("B50").Value = 'BabyDogeCoin'
$ws.Range("C50").Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.00000000117'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -1.02%  '
$ws.Range("B51").Value = 'TheSandbox'
$ws.Range("C51").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.4030'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -0.82%  '
